# Add Original Color option for pptx
# - Change slide master background to a solid blue fill (instead of the
#   scheme bg1 reference)
# - Change all slide title/content text color from black to white
# - Append two new Q&A slides (Question 5 / Answer) matching the style
#   of the existing Q&A slides

$p = $ppt.ActivePresentation

# 1. Slide master background: solid blue fill instead of scheme bg1 ref
$master = $p.SlideMaster
$master.Background.Fill.Solid()
$master.Background.Fill.ForeColor.RGB = 16711680

# 2. Add the two new Q&A slides (Question 5 / Answer) at the end, by
#    duplicating the last existing Question/Answer pair so the new
#    slides inherit the same placeholder layout & run formatting
#    (font sizes) without forcing an autofit recalculation.
$srcQuestion = $p.Slides.Item(8)
$srcAnswer = $p.Slides.Item(9)

$newQuestion = $srcQuestion.Duplicate()
$newQuestion.MoveTo($p.Slides.Count)
$newQuestion.Shapes.Item(1).TextFrame.TextRange.Text = "Question 5"
$newQuestion.Shapes.Item(2).TextFrame.TextRange.Text = "What is Market Failure?"

$newAnswer = $srcAnswer.Duplicate()
$newAnswer.MoveTo($p.Slides.Count)
$newAnswer.Shapes.Item(1).TextFrame.TextRange.Text = "Answer"
$newAnswer.Shapes.Item(2).TextFrame.TextRange.Text = "Market failure is the economic situation defined by an inefficient distribution of goods and services in the free market."

# 3. Recolor every slide's text runs (black -> white), including the
#    two slides that were just appended.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Length -gt 0) {
                $tr.Font.Color.RGB = 16777215
            }
        }
    }
}
